# Change title slide location/date from "Lincoln, NE – November 19, 2016"
# to "Madison, WI – April 8, 2017"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$para2 = $tr.Paragraphs(2)
$para2.Text = "Madison, WI – April 8"
$null = $para2.InsertAfter(", 2017")
